# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1815
$ws1.Range("F10").Value = 3587
$ws1.Range("F16").Value = 622
$ws1.Range("F17").Value = 116
$ws1.Range("F25").Value = 2826
$ws1.Range("F26").Value = 5290
$ws1.Range("F30").Value = 3107
$ws1.Range("F31").Value = 301
$ws1.Range("F32").Value = 2295
$ws1.Range("F34").Value = 495
$ws1.Range("F40").Value = 473
$ws1.Range("F41").Value = 819
$ws1.Range("F44").Value = 0
$ws1.Range("F46").Value = 505

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1815
$ws4.Range("F10").Value = 3587
$ws4.Range("F17").Value = 622
$ws4.Range("F18").Value = 116
$ws4.Range("F26").Value = 2826
$ws4.Range("F27").Value = 5290
$ws4.Range("F31").Value = 3107
$ws4.Range("F32").Value = 301
$ws4.Range("F33").Value = 2295
$ws4.Range("F35").Value = 495
$ws4.Range("F41").Value = 473
$ws4.Range("F42").Value = 819
$ws4.Range("F47").Value = 505
